$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.243534922599792
$ws.Range("B1").Value = 2.194370985031128
$ws.Range("C1").Value = 4.263675212860107
$ws.Range("D1").Value = 3.038837909698486
$ws.Range("E1").Value = 1.055611133575439
